# Natmi following Dr Hou advice
# Ligand-/receptor-expressing cell counts for L1cam-Itgav go from 1 to 3
# detected cells in every sending/target cluster combination, which
# changes the dependent average/total expression and specificity values
# for rows 2-10 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-D (clusters, symbols) and F, L (detection rate) are untouched.
$rowData = @{
    2  = @{ E=3; G=22.59487733333333;  H=67.784632;           I=0.7395019553569895;  J=0.7395019553569895;  K=3; M=21.09934133333334; N=63.29802400000001; O=0.2917236204149438; P=0.2917236204149438; Q=476.7370292407965;  R=4290.633263167168;  S=0.2157301877206711;   T=0.2157301877206711 }
    3  = @{ E=3; G=22.59487733333333;  H=67.784632;           I=0.7395019553569895;  J=0.7395019553569895;  K=3; M=35.81943766666667; N=107.458313;          O=0.4952465516465762; P=0.4952465516465762; Q=809.3358002273128;  R=7284.022202045817;  S=0.3662357933264493;   T=0.3662357933264494 }
    4  = @{ E=3; G=22.59487733333333;  H=67.784632;           I=0.7395019553569895;  J=0.7395019553569895;  K=3; M=15.40769666666667; N=46.22309;            O=0.2130298279384801; P=0.2130298279384801; Q=348.1350161725422;  R=3133.21514555288;   S=0.157535974309869;    T=0.1575359743098691 }
    5  = @{ E=3; G=0.3045986666666667; H=0.9137960000000001; I=0.00996913177602551; J=0.00996913177602551; K=3; M=21.09934133333334; N=63.29802400000001; O=0.2917236204149438; P=0.2917236204149438; Q=6.426831237678223;  R=57.841481139104;    S=0.00290823121409582;  T=0.002908231214095821 }
    6  = @{ E=3; G=0.3045986666666667; H=0.9137960000000001; I=0.00996913177602551; J=0.00996913177602551; K=3; M=35.81943766666667; N=107.458313;          O=0.4952465516465762; P=0.4952465516465762; Q=10.91055295401645;  R=98.19497658614802;  S=0.004937178134986941; T=0.004937178134986942 }
    7  = @{ E=3; G=0.3045986666666667; H=0.9137960000000001; I=0.00996913177602551; J=0.00996913177602551; K=3; M=15.40769666666667; N=46.22309;            O=0.2130298279384801; P=0.2130298279384801; Q=4.693163861071111;  R=42.23847474964;     S=0.002123722426942749; T=0.002123722426942749 }
    8  = @{ E=3; G=7.654706;           H=22.964118;            I=0.2505289128669849;  J=0.2505289128669849;  K=3; M=21.09934133333334; N=63.29802400000001; O=0.2917236204149438; P=0.2917236204149438; Q=161.5092547003147;  R=1453.583292302832;  S=0.07308520148017683;  T=0.07308520148017685 }
    9  = @{ E=3; G=7.654706;           H=22.964118;            I=0.2505289128669849;  J=0.2505289128669849;  K=3; M=35.81943766666667; N=107.458313;          O=0.4952465516465762; P=0.4952465516465762; Q=274.1872644236593;  R=2467.685379812934;  S=0.1240735801851398;   T=0.1240735801851398 }
    10 = @{ E=3; G=7.654706;           H=22.964118;            I=0.2505289128669849;  J=0.2505289128669849;  K=3; M=15.40769666666667; N=46.22309;            O=0.2130298279384801; P=0.2130298279384801; Q=117.9413881205133;  R=1061.47249308462;   S=0.05337013120166827;  T=0.05337013120166827 }
}

foreach ($r in $rowData.Keys) {
    $cols = $rowData[$r]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$r").Value = $cols[$col]
    }
}
